# Generate Report for Handoff
# Updates the "a16ed3cc-7b04-4cfd-92e1-6a68ddf62816.md" file's handoff / xliff
# generation timestamps across the Overview, zh-cn and de-de sheets to reflect
# a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# Column G = "Latest HO Xliff Generate Date" for the a16ed3cc row (row 6)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-19 00:40:39"

# --- zh-cn sheet ------------------------------------------------------------
# Column H = "Latest Handoff Datetime" for the a16ed3cc row (row 6)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-19 00:40:33"

# --- de-de sheet ------------------------------------------------------------
# Column H = "Latest Handoff Datetime" for the a16ed3cc row (row 6)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-19 00:40:39"
